$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(242, 1).Value = 46021
$ws.Cells.Item(242, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(242, 3).Value = 450.23099999999999
$ws.Cells.Item(242, 4).Value = 927.96400000000017
$ws.Cells.Item(242, 5).Value = 298.94199999999995
$ws.Cells.Item(242, 6).Value = 230.80199999999996
$ws.Cells.Item(242, 7).Value = 183.19699999999997
$ws.Cells.Item(242, 8).Value = 443.97500000000002
$ws.Cells.Item(242, 9).Value = 649.09099999999989
$ws.Cells.Item(242, 10).Value = 191.62099999999998
$ws.Cells.Item(242, 11).Value = 185.75500000000002
$ws.Cells.Item(242, 12).Value = 166.01000000000002
$ws.Cells.Item(242, 13).Value = 131.47999999999999
$ws.Cells.Item(242, 14).Value = 218.66
$ws.Cells.Item(242, 15).Value = 642.95900000000006
$ws.Cells.Item(242, 16).Value = 1454.6390000000004
$ws.Cells.Item(242, 17).Value = 374.44
$ws.Cells.Item(242, 18).Value = 501.95499999999998
$ws.Cells.Item(242, 19).Value = 374.92399999999998
$ws.Cells.Item(242, 20).Value = 253.92200000000003
$ws.Cells.Item(242, 21).Value = 71.819999999999993
$ws.Cells.Item(242, 22).Value = 45.134
$ws.Cells.Item(242, 23).Value = 22.24
$ws.Cells.Item(242, 24).Value = 42.06
$ws.Cells.Item(242, 25).Value = 96.78
$ws.Cells.Item(242, 26).Value = 88.86

$ws.Cells.Item(243, 1).Value = 46021
$ws.Cells.Item(243, 2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(243, 3).Value = 282.90299999999996
$ws.Cells.Item(243, 4).Value = 736.18700000000001
$ws.Cells.Item(243, 5).Value = 383.12300000000005
$ws.Cells.Item(243, 6).Value = 58.122
$ws.Cells.Item(243, 7).Value = 59.213999999999999
$ws.Cells.Item(243, 8).Value = 252.97499999999997
$ws.Cells.Item(243, 9).Value = 68.88300000000001
$ws.Cells.Item(243, 10).Value = 107.55600000000001
$ws.Cells.Item(243, 11).Value = 129.74200000000002
$ws.Cells.Item(243, 12).Value = 218.416
$ws.Cells.Item(243, 13).Value = 83.518000000000001
$ws.Cells.Item(243, 14).Value = 223.45000000000002
$ws.Cells.Item(243, 15).Value = 471.73199999999997
$ws.Cells.Item(243, 16).Value = 844.30199999999991
$ws.Cells.Item(243, 17).Value = 325.98500000000001
$ws.Cells.Item(243, 18).Value = 201.93399999999997
$ws.Cells.Item(243, 19).Value = 112.82300000000001
$ws.Cells.Item(243, 20).Value = 164.58600000000001
$ws.Cells.Item(243, 21).Value = 168.422
$ws.Cells.Item(243, 22).Value = 74.590999999999994
$ws.Cells.Item(243, 23).Value = 114.898
$ws.Cells.Item(243, 24).Value = 127.14100000000001
$ws.Cells.Item(243, 25).Value = 0
$ws.Cells.Item(243, 26).Value = 0

$ws.Cells.Item(244, 1).Value = 46022
$ws.Cells.Item(244, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(244, 3).Value = 571.67999999999995
$ws.Cells.Item(244, 4).Value = 517.84
$ws.Cells.Item(244, 5).Value = 508.50700000000001
$ws.Cells.Item(244, 6).Value = 269.06599999999997
$ws.Cells.Item(244, 7).Value = 147.727
$ws.Cells.Item(244, 8).Value = 357.11399999999992
$ws.Cells.Item(244, 9).Value = 298.41100000000006
$ws.Cells.Item(244, 10).Value = 173.59900000000002
$ws.Cells.Item(244, 11).Value = 116.98299999999999
$ws.Cells.Item(244, 12).Value = 92.13
$ws.Cells.Item(244, 13).Value = 42.03
$ws.Cells.Item(244, 14).Value = 95.7
$ws.Cells.Item(244, 15).Value = 222.68000000000004
$ws.Cells.Item(244, 16).Value = 246.5
$ws.Cells.Item(244, 17).Value = 145.5
$ws.Cells.Item(244, 18).Value = 66.14
$ws.Cells.Item(244, 19).Value = 36.94
$ws.Cells.Item(244, 20).Value = 84.539999999999992
$ws.Cells.Item(244, 21).Value = 112.18
$ws.Cells.Item(244, 22).Value = 6.46
$ws.Cells.Item(244, 23).Value = 0
$ws.Cells.Item(244, 24).Value = 0
$ws.Cells.Item(244, 25).Value = 0.74
$ws.Cells.Item(244, 26).Value = 81.84

$ws.Cells.Item(245, 1).Value = 46022
$ws.Cells.Item(245, 2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(245, 3).Value = 413.58699999999999
$ws.Cells.Item(245, 4).Value = 492.01100000000002
$ws.Cells.Item(245, 5).Value = 305.30899999999997
$ws.Cells.Item(245, 6).Value = 61.816999999999993
$ws.Cells.Item(245, 7).Value = 258.73500000000001
$ws.Cells.Item(245, 8).Value = 145.22800000000001
$ws.Cells.Item(245, 9).Value = 237.37299999999999
$ws.Cells.Item(245, 10).Value = 283.07
$ws.Cells.Item(245, 11).Value = 190.50800000000001
$ws.Cells.Item(245, 12).Value = 188.358
$ws.Cells.Item(245, 13).Value = 115.67299999999999
$ws.Cells.Item(245, 14).Value = 255.02299999999997
$ws.Cells.Item(245, 15).Value = 586.11899999999991
$ws.Cells.Item(245, 16).Value = 1016.4630000000001
$ws.Cells.Item(245, 17).Value = 201.6
$ws.Cells.Item(245, 18).Value = 209.24600000000001
$ws.Cells.Item(245, 19).Value = 184.10299999999998
$ws.Cells.Item(245, 20).Value = 127.33000000000001
$ws.Cells.Item(245, 21).Value = 48.68
$ws.Cells.Item(245, 22).Value = 12.281000000000001
$ws.Cells.Item(245, 23).Value = 46.923999999999999
$ws.Cells.Item(245, 24).Value = 151.613
$ws.Cells.Item(245, 25).Value = 34.367999999999995
$ws.Cells.Item(245, 26).Value = 62.401000000000003

$ws.Cells.Item(246, 1).Value = 46023
$ws.Cells.Item(246, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(246, 3).Value = 652.52
$ws.Cells.Item(246, 4).Value = 1363.22
$ws.Cells.Item(246, 5).Value = 953.31
$ws.Cells.Item(246, 6).Value = 756.58
$ws.Cells.Item(246, 7).Value = 477.43
$ws.Cells.Item(246, 8).Value = 688.85
$ws.Cells.Item(246, 9).Value = 827.54
$ws.Cells.Item(246, 10).Value = 171.22
$ws.Cells.Item(246, 11).Value = 257.81
$ws.Cells.Item(246, 12).Value = 365.46
$ws.Cells.Item(246, 13).Value = 203.21
$ws.Cells.Item(246, 14).Value = 263.18
$ws.Cells.Item(246, 15).Value = 713.86
$ws.Cells.Item(246, 16).Value = 1592.91
$ws.Cells.Item(246, 17).Value = 1282.8900000000001
$ws.Cells.Item(246, 18).Value = 494.85
$ws.Cells.Item(246, 19).Value = 600.5
$ws.Cells.Item(246, 20).Value = 356.94
$ws.Cells.Item(246, 21).Value = 261.10000000000002
$ws.Cells.Item(246, 22).Value = 125.25
$ws.Cells.Item(246, 23).Value = 206.14
$ws.Cells.Item(246, 24).Value = 98.98
$ws.Cells.Item(246, 25).Value = 80.16
$ws.Cells.Item(246, 26).Value = 25.16

$ws.Cells.Item(247, 1).Value = 46023
$ws.Cells.Item(247, 2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(247, 3).Value = 417.95699999999999
$ws.Cells.Item(247, 4).Value = 743.40100000000007
$ws.Cells.Item(247, 5).Value = 200.33100000000002
$ws.Cells.Item(247, 6).Value = 32.914000000000001
$ws.Cells.Item(247, 7).Value = 84.140999999999991
$ws.Cells.Item(247, 8).Value = 138.471
$ws.Cells.Item(247, 9).Value = 189.50700000000001
$ws.Cells.Item(247, 10).Value = 0
$ws.Cells.Item(247, 11).Value = 38.875999999999998
$ws.Cells.Item(247, 12).Value = 66.670999999999992
$ws.Cells.Item(247, 13).Value = 88.406000000000006
$ws.Cells.Item(247, 14).Value = 197.572
$ws.Cells.Item(247, 15).Value = 408.95600000000002
$ws.Cells.Item(247, 16).Value = 546.70999999999992
$ws.Cells.Item(247, 17).Value = 418.589
$ws.Cells.Item(247, 18).Value = 303.09699999999998
$ws.Cells.Item(247, 19).Value = 202.709
$ws.Cells.Item(247, 20).Value = 153.47299999999998
$ws.Cells.Item(247, 21).Value = 72.492000000000004
$ws.Cells.Item(247, 22).Value = 142.76599999999999
$ws.Cells.Item(247, 23).Value = 52.707000000000001
$ws.Cells.Item(247, 24).Value = 64.385999999999996
$ws.Cells.Item(247, 25).Value = 43.203000000000003
$ws.Cells.Item(247, 26).Value = 0

$ws.Cells.Item(248, 1).Value = 46024
$ws.Cells.Item(248, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(248, 3).Value = 639.6
$ws.Cells.Item(248, 4).Value = 1638.38
$ws.Cells.Item(248, 5).Value = 831.8
$ws.Cells.Item(248, 6).Value = 649.44000000000005
$ws.Cells.Item(248, 7).Value = 344.77
$ws.Cells.Item(248, 8).Value = 728.9
$ws.Cells.Item(248, 9).Value = 429.98
$ws.Cells.Item(248, 10).Value = 234.03
$ws.Cells.Item(248, 11).Value = 109.93
$ws.Cells.Item(248, 12).Value = 200.66
$ws.Cells.Item(248, 13).Value = 190.5
$ws.Cells.Item(248, 14).Value = 455.11
$ws.Cells.Item(248, 15).Value = 796.08
$ws.Cells.Item(248, 16).Value = 1764.99
$ws.Cells.Item(248, 17).Value = 623.33000000000004
$ws.Cells.Item(248, 18).Value = 502.74
$ws.Cells.Item(248, 19).Value = 266.98
$ws.Cells.Item(248, 20).Value = 307.99
$ws.Cells.Item(248, 21).Value = 116.73
$ws.Cells.Item(248, 22).Value = 82.02
$ws.Cells.Item(248, 23).Value = 62.17
$ws.Cells.Item(248, 24).Value = 143.38
$ws.Cells.Item(248, 25).Value = 95.6
$ws.Cells.Item(248, 26).Value = 24.82

$ws.Cells.Item(249, 1).Value = 46024
$ws.Cells.Item(249, 2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(249, 3).Value = 471.34400000000005
$ws.Cells.Item(249, 4).Value = 751.96
$ws.Cells.Item(249, 5).Value = 207.03900000000002
$ws.Cells.Item(249, 6).Value = 65.289999999999992
$ws.Cells.Item(249, 7).Value = 72.733999999999995
$ws.Cells.Item(249, 8).Value = 169.34399999999999
$ws.Cells.Item(249, 9).Value = 222.554
$ws.Cells.Item(249, 10).Value = 8.93
$ws.Cells.Item(249, 11).Value = 180.38400000000001
$ws.Cells.Item(249, 12).Value = 152.69
$ws.Cells.Item(249, 13).Value = 22.398
$ws.Cells.Item(249, 14).Value = 90.205000000000013
$ws.Cells.Item(249, 15).Value = 509.99599999999998
$ws.Cells.Item(249, 16).Value = 607.63900000000001
$ws.Cells.Item(249, 17).Value = 779.19699999999989
$ws.Cells.Item(249, 18).Value = 119.79500000000002
$ws.Cells.Item(249, 19).Value = 292.61500000000001
$ws.Cells.Item(249, 20).Value = 111.25999999999999
$ws.Cells.Item(249, 21).Value = 120.10999999999999
$ws.Cells.Item(249, 22).Value = 59.039000000000001
$ws.Cells.Item(249, 23).Value = 0
$ws.Cells.Item(249, 24).Value = 15.744999999999999
$ws.Cells.Item(249, 25).Value = 10.367000000000001
$ws.Cells.Item(249, 26).Value = 87.905000000000001

$ws.Cells.Item(250, 1).Value = 46025
$ws.Cells.Item(250, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(250, 3).Value = 815.64
$ws.Cells.Item(250, 4).Value = 1626.81
$ws.Cells.Item(250, 5).Value = 1042.0899999999999
$ws.Cells.Item(250, 6).Value = 665.47
$ws.Cells.Item(250, 7).Value = 349.92
$ws.Cells.Item(250, 8).Value = 773.27
$ws.Cells.Item(250, 9).Value = 377.11
$ws.Cells.Item(250, 10).Value = 212.04
$ws.Cells.Item(250, 11).Value = 124.52
$ws.Cells.Item(250, 12).Value = 171.29
$ws.Cells.Item(250, 13).Value = 215.89
$ws.Cells.Item(250, 14).Value = 303.61
$ws.Cells.Item(250, 15).Value = 769.46
$ws.Cells.Item(250, 16).Value = 1793.19
$ws.Cells.Item(250, 17).Value = 971.51
$ws.Cells.Item(250, 18).Value = 906.47
$ws.Cells.Item(250, 19).Value = 338.83
$ws.Cells.Item(250, 20).Value = 199.85
$ws.Cells.Item(250, 21).Value = 171.14
$ws.Cells.Item(250, 22).Value = 131.02000000000001
$ws.Cells.Item(250, 23).Value = 105.91
$ws.Cells.Item(250, 24).Value = 230.6
$ws.Cells.Item(250, 25).Value = 151.1
$ws.Cells.Item(250, 26).Value = 15.66

$ws.Cells.Item(251, 1).Value = 46025
$ws.Cells.Item(251, 2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(251, 3).Value = 259.815
$ws.Cells.Item(251, 4).Value = 689.68000000000006
$ws.Cells.Item(251, 5).Value = 205.27500000000001
$ws.Cells.Item(251, 6).Value = 95.649000000000001
$ws.Cells.Item(251, 7).Value = 158.86699999999999
$ws.Cells.Item(251, 8).Value = 130.946
$ws.Cells.Item(251, 9).Value = 322.161
$ws.Cells.Item(251, 10).Value = 18.368000000000002
$ws.Cells.Item(251, 11).Value = 81.73599999999999
$ws.Cells.Item(251, 12).Value = 189.66500000000002
$ws.Cells.Item(251, 13).Value = 61.050999999999995
$ws.Cells.Item(251, 14).Value = 73.91
$ws.Cells.Item(251, 15).Value = 539.70500000000004
$ws.Cells.Item(251, 16).Value = 558.1629999999999
$ws.Cells.Item(251, 17).Value = 276.89500000000004
$ws.Cells.Item(251, 18).Value = 220.51299999999998
$ws.Cells.Item(251, 19).Value = 146.251
$ws.Cells.Item(251, 20).Value = 151.124
$ws.Cells.Item(251, 21).Value = 63.411000000000001
$ws.Cells.Item(251, 22).Value = 154.352
$ws.Cells.Item(251, 23).Value = 59.530999999999999
$ws.Cells.Item(251, 24).Value = 0
$ws.Cells.Item(251, 25).Value = 8.6720000000000006
$ws.Cells.Item(251, 26).Value = 5.992

$ws.Cells.Item(252, 1).Value = 46026
$ws.Cells.Item(252, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(252, 3).Value = 699.63
$ws.Cells.Item(252, 4).Value = 2073.04
$ws.Cells.Item(252, 5).Value = 1124.81
$ws.Cells.Item(252, 6).Value = 352.78
$ws.Cells.Item(252, 7).Value = 325.92
$ws.Cells.Item(252, 8).Value = 658.93
$ws.Cells.Item(252, 9).Value = 342.34
$ws.Cells.Item(252, 10).Value = 166.88
$ws.Cells.Item(252, 11).Value = 270.05
$ws.Cells.Item(252, 12).Value = 170.57
$ws.Cells.Item(252, 13).Value = 249.79
$ws.Cells.Item(252, 14).Value = 400.33
$ws.Cells.Item(252, 15).Value = 1033.0999999999999
$ws.Cells.Item(252, 16).Value = 2882.27
$ws.Cells.Item(252, 17).Value = 1248.1600000000001
$ws.Cells.Item(252, 18).Value = 444.83
$ws.Cells.Item(252, 19).Value = 361.24
$ws.Cells.Item(252, 20).Value = 206.54
$ws.Cells.Item(252, 21).Value = 111.81
$ws.Cells.Item(252, 22).Value = 31.57
$ws.Cells.Item(252, 23).Value = 111.66
$ws.Cells.Item(252, 24).Value = 119.22
$ws.Cells.Item(252, 25).Value = 64.72
$ws.Cells.Item(252, 26).Value = 138.63

$ws.Cells.Item(253, 1).Value = 46026
$ws.Cells.Item(253, 2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(253, 3).Value = 487.99899999999997
$ws.Cells.Item(253, 4).Value = 760.447
$ws.Cells.Item(253, 5).Value = 204.73499999999999
$ws.Cells.Item(253, 6).Value = 101.899
$ws.Cells.Item(253, 7).Value = 37.055999999999997
$ws.Cells.Item(253, 8).Value = 104.99600000000001
$ws.Cells.Item(253, 9).Value = 205.154
$ws.Cells.Item(253, 10).Value = 32.162999999999997
$ws.Cells.Item(253, 11).Value = 102.43599999999999
$ws.Cells.Item(253, 12).Value = 202.37299999999999
$ws.Cells.Item(253, 13).Value = 208.58300000000003
$ws.Cells.Item(253, 14).Value = 255.15
$ws.Cells.Item(253, 15).Value = 822.04199999999992
$ws.Cells.Item(253, 16).Value = 745.02599999999995
$ws.Cells.Item(253, 17).Value = 395.20800000000008
$ws.Cells.Item(253, 18).Value = 695.07499999999982
$ws.Cells.Item(253, 19).Value = 653.93399999999997
$ws.Cells.Item(253, 20).Value = 618.1629999999999
$ws.Cells.Item(253, 21).Value = 89.544999999999987
$ws.Cells.Item(253, 22).Value = 69.343000000000004
$ws.Cells.Item(253, 23).Value = 111.133
$ws.Cells.Item(253, 24).Value = 18.167999999999999
$ws.Cells.Item(253, 25).Value = 31.426000000000002
$ws.Cells.Item(253, 26).Value = 0

$ws.Range("G257").Select()
Write-Host "done"